{"js": "// Update the 25 \"three-digit \u00f7 one-digit\" practice answers in the table.\n// Each old value is unique in the document, so a direct search+replace\n// (by exact, case-sensitive match) for each pair is safe and unambiguous.\nconst replacements = [\n  [\"702\u00f76=117, 0\", \"441\u00f79=49, 0\"],\n  [\"332\u00f72=166, 0\", \"270\u00f72=135, 0\"],\n  [\"924\u00f74=231, 0\", \"768\u00f75=153, 3\"],\n  [\"957\u00f78=119, 5\", \"367\u00f77=52, 3\"],\n  [\"230\u00f76=38, 2\", \"742\u00f76=123, 4\"],\n  [\"529\u00f72=264, 1\", \"270\u00f75=54, 0\"],\n  [\"685\u00f74=171, 1\", \"572\u00f74=143, 0\"],\n  [\"352\u00f78=44, 0\", \"239\u00f78=29, 7\"],\n  [\"370\u00f74=92, 2\", \"325\u00f78=40, 5\"],\n  [\"238\u00f73=79, 1\", \"342\u00f79=38, 0\"],\n  [\"746\u00f72=373, 0\", \"783\u00f73=261, 0\"],\n  [\"257\u00f78=32, 1\", \"746\u00f74=186, 2\"],\n  [\"107\u00f76=17, 5\", \"174\u00f76=29, 0\"],\n  [\"365\u00f76=60, 5\", \"711\u00f73=237, 0\"],\n  [\"318\u00f77=45, 3\", \"532\u00f72=266, 0\"],\n  [\"985\u00f72=492, 1\", \"367\u00f73=122, 1\"],\n  [\"573\u00f74=143, 1\", \"233\u00f74=58, 1\"],\n  [\"464\u00f75=92, 4\", \"660\u00f78=82, 4\"],\n  [\"521\u00f78=65, 1\", \"851\u00f77=121, 4\"],\n  [\"634\u00f78=79, 2\", \"242\u00f76=40, 2\"],\n  [\"228\u00f76=38, 0\", \"218\u00f78=27, 2\"],\n  [\"354\u00f72=177, 0\", \"932\u00f77=133, 1\"],\n  [\"573\u00f78=71, 5\", \"929\u00f77=132, 5\"],\n  [\"399\u00f72=199, 1\", \"129\u00f74=32, 1\"],\n  [\"782\u00f78=97, 6\", \"446\u00f74=111, 2\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the 25 \"three-digit \u00f7 one-digit\" practice answers in the table.\n# Each old value is unique in the document, so Find/Replace (exact match,\n# whole document, one hit each) is safe and unambiguous.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Old = \"702\u00f76=117, 0\"; New = \"441\u00f79=49, 0\" },\n    @{ Old = \"332\u00f72=166, 0\"; New = \"270\u00f72=135, 0\" },\n    @{ Old = \"924\u00f74=231, 0\"; New = \"768\u00f75=153, 3\" },\n    @{ Old = \"957\u00f78=119, 5\"; New = \"367\u00f77=52, 3\" },\n    @{ Old = \"230\u00f76=38, 2\"; New = \"742\u00f76=123, 4\" },\n    @{ Old = \"529\u00f72=264, 1\"; New = \"270\u00f75=54, 0\" },\n    @{ Old = \"685\u00f74=171, 1\"; New = \"572\u00f74=143, 0\" },\n    @{ Old = \"352\u00f78=44, 0\"; New = \"239\u00f78=29, 7\" },\n    @{ Old = \"370\u00f74=92, 2\"; New = \"325\u00f78=40, 5\" },\n    @{ Old = \"238\u00f73=79, 1\"; New = \"342\u00f79=38, 0\" },\n    @{ Old = \"746\u00f72=373, 0\"; New = \"783\u00f73=261, 0\" },\n    @{ Old = \"257\u00f78=32, 1\"; New = \"746\u00f74=186, 2\" },\n    @{ Old = \"107\u00f76=17, 5\"; New = \"174\u00f76=29, 0\" },\n    @{ Old = \"365\u00f76=60, 5\"; New = \"711\u00f73=237, 0\" },\n    @{ Old = \"318\u00f77=45, 3\"; New = \"532\u00f72=266, 0\" },\n    @{ Old = \"985\u00f72=492, 1\"; New = \"367\u00f73=122, 1\" },\n    @{ Old = \"573\u00f74=143, 1\"; New = \"233\u00f74=58, 1\" },\n    @{ Old = \"464\u00f75=92, 4\"; New = \"660\u00f78=82, 4\" },\n    @{ Old = \"521\u00f78=65, 1\"; New = \"851\u00f77=121, 4\" },\n    @{ Old = \"634\u00f78=79, 2\"; New = \"242\u00f76=40, 2\" },\n    @{ Old = \"228\u00f76=38, 0\"; New = \"218\u00f78=27, 2\" },\n    @{ Old = \"354\u00f72=177, 0\"; New = \"932\u00f77=133, 1\" },\n    @{ Old = \"573\u00f78=71, 5\"; New = \"929\u00f77=132, 5\" },\n    @{ Old = \"399\u00f72=199, 1\"; New = \"129\u00f74=32, 1\" },\n    @{ Old = \"782\u00f78=97, 6\"; New = \"446\u00f74=111, 2\" }\n)\n\nforeach ($pair in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $pair.Old\n    $find.Replacement.Text = $pair.New\n    $find.Execute($pair.Old, $true, $false, $false, $false, $false, $true, 1, $false, $pair.New, 2)\n}\n"}
